# =====================================================================
# Adds a new "96 Total" column (CG) to the two_digit_data worksheet,
# containing the row totals for each monthly observation, and updates
# the sheet selection to reflect the newly-active cell near the new
# column (matching the author's edit).
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("two_digit_data")

# ---------------------------------------------------------------
# 1. Header cell CG1: "96 Total", formatted like the other
#    category headers (copy the format from CF1, the previous
#    last header cell).
# ---------------------------------------------------------------
$ws.Range("CF1").Copy()
$ws.Range("CG1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("CG1").Value = "96 Total"

# ---------------------------------------------------------------
# 2. Data cells CG2:CG152: totals, formatted with the same
#    number format used elsewhere in the sheet (numFmtId 164)
#    but on the small Arial 8 font (fontId 20) already used by
#    other helper cells (e.g. BC153) - this reproduces the new
#    style added to cellXfs.
# ---------------------------------------------------------------
$ws.Range("BC153").Copy()
$ws.Range("CG2:CG152").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("CG2:CG152").NumberFormat = "0.0;\-0.0;0.0;@"

$totals = @(
    6576.0618182999997, 6595.6641925000004, 6657.9327706000004, 6675.5378567999996, 6881.9441951999997, 6861.8554940000004, 6978.9699757999997, 6918.5250524000003,
    7034.1639395000002, 7022.2527528999999, 7103.5234570000002, 7092.2573214000004, 7182.1613009000002, 7247.3043494000003, 7354.8338311999996, 7353.2863914999998,
    7508.606914, 7566.4884234000001, 7720.7895072000001, 7715.2730129000001, 7862.3568226999996, 7802.9902621000001, 7908.4636763999997, 7807.9828870000001,
    7826.9227229999997, 7708.7805939, 7703.7164486000001, 7620.7032112999996, 7614.9417599999997, 7597.1754529999998, 7631.4000545999997, 7613.3991011999997,
    7587.4465005000002, 7530.9593086000004, 7620.5067890999999, 7589.4031314000003, 7744.9260545999996, 7738.5979183999998, 7857.9087863000004, 7861.6938744999998,
    8003.2533721999998, 8049.8759320999998, 8165.2016004999996, 8165.9187173999999, 8300.2621098, 8246.6034779000001, 8297.1254308000007, 8265.8648565000003,
    8314.3698521999995, 8291.3875501999992, 8330.1864401999992, 8250.3060526000008, 8462.0847288000004, 8412.9965838999997, 8482.5400183999991, 8455.2334644000002,
    8554.2948639000006, 8556.4741814999998, 8606.5201228000005, 8603.6235892000004, 8763.6674767999993, 8754.4042205000005, 8864.9997810999994, 8911.5878905000009,
    8912.7194727000006, 8920.9354323999996, 9004.6291142999999, 8971.8953963000004, 9055.6440251000004, 9092.3099364999998, 9134.6730812999995, 9150.2861487999999,
    9286.7409043999996, 9372.5993873000007, 9380.0778465999992, 9326.9510515999991, 9427.7396370000006, 9434.4538051, 9533.4733240000005, 9462.0865529000002,
    9652.9227023999993, 9744.6145959999994, 9850.7550761000002, 9870.3355114000005, 9944.9981674999999, 9961.1373913999996, 10082.055874600001, 10105.7732289,
    10201.645047, 10276.9833392, 10412.0073436, 10406.3975414, 10535.100047800001, 10622.906644500001, 10690.3063561, 10710.441137600001,
    10757.3446362, 10776.7288095, 10786.7987139, 10707.257287500001, 10829.231529500001, 10889.154493399999, 10957.4266088, 10973.8229331,
    11144.5452391, 11146.1760954, 11196.181066900001, 11127.2154885, 11239.4184518, 11237.735913099999, 11393.152977899999, 11264.1168555,
    11378.3335893, 11425.3118584, 11491.664766399999, 11361.409787299999, 11445.9000442, 11464.564434100001, 11553.348763100001, 11572.7635314,
    11581.725704, 11724.132317, 11781.440931200001, 11702.631254899999, 11935.266320000001, 11975.894034499999, 11999.9494196, 11904.7089265,
    12054.006582399999, 12103.333139, 12256.6339114, 12245.332033999999, 12453.6309795, 12528.1846948, 12572.697362000001, 12547.5377006,
    12741.586741900001, 12802.157754600001, 12921.060775899999, 12858.8226376, 12996.176518099999, 13053.282018100001, 12192.5802874, 12562.7918001,
    12882.160213700001, 13037.041243600001, 13204.3429452, 12955.3726811, 13273.977877699999, 13482.1191734, 13589.1047423
)

$data = New-Object 'object[,]' 151,1
for ($i = 0; $i -lt 151; $i++) {
    $data[$i, 0] = $totals[$i]
}
$ws.Range("CG2:CG152").Value = $data

# ---------------------------------------------------------------
# 3. View state: scroll towards the new column and move the
#    active selection to CG2 (the first new total cell), matching
#    the author's saved selection/scroll position.
# ---------------------------------------------------------------
$ws.Range("BV1").Select()
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 74   # column BV
$aw.ScrollRow = 1
$ws.Range("CG2").Select()
